# Apply cryptos.xlsx data refresh (Mon May 13 11:54:40 UTC 2024 GitHub Actions run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-coerced to a Number by Excel
# (plain decimals like "81.90" or "0.0000234") are written via a temporary
# Text number-format so the literal string survives, then the format is reset
# back to Normal/General so no stray style is left on the cell.
$ws.Range("D2").Value = '62.689.96'
$ws.Range("E2").Value = '  +2.58%  '
$ws.Range("D3").Value = '2.963.38'
$ws.Range("E3").Value = '  +1.28%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '595.35'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.55%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '145.85'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.55%  '
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("D8").Value = '2.959.82'
$ws.Range("E8").Value = '  +1.17%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.506'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.36%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.35'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +4.73%  '
$ws.Range("E11").Value = '  -0.12%  '
$ws.Range("E12").Value = '  +1.73%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000234'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +3.86%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '33.37'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.18%  '
$ws.Range("E15").Value = '  -0.22%  '
$ws.Range("D16").Value = '3.450.27'
$ws.Range("E16").Value = '  +1.15%  '
$ws.Range("D17").Value = '62.564.25'
$ws.Range("E17").Value = '  +2.40%  '
$ws.Range("B18").Value = 'WrappedEther'
$ws.Range("C18").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D18").Value = '2.989.09'
$ws.Range("E18").Value = '  +2.03%  '
$ws.Range("B19").Value = 'Polkadot'
$ws.Range("C19").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.71'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.28%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '442.01'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.44%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.42'
$ws.Range("D21").Style = "Normal"
$ws.Range("E22").Value = '  -1.06%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.09'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.36%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '81.90'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.58%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '11.25'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.62%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.97'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.72%  '
$ws.Range("E27").Value = '  -3.51%  '
$ws.Range("E28").Value = '  -0.01%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.60'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.10%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.10'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.65%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.12'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -5.94%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '26.67'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.24%  '
$ws.Range("E33").Value = '  -2.24%  '
$ws.Range("E34").Value = '  -0.09%  '
$ws.Range("D35").Value = '0.0₃0881'
$ws.Range("E35").Value = '  +1.74%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.995'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.30%  '
$ws.Range("E37").Value = '  +0.04%  '
$ws.Range("B38").Value = 'Stacks'
$ws.Range("C38").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.05'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.56%  '
$ws.Range("B39").Value = 'OKB'
$ws.Range("C39").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '49.63'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.21%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.93'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.13%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '8.57'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.35%  '
$ws.Range("E42").Value = '  -4.86%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.282'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.37%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '39.07'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -6.87%  '
$ws.Range("D45").Value = '2.722.88'
$ws.Range("E45").Value = '  +1.36%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '136.20'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.98%  '
$ws.Range("E47").Value = '  -1.56%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '364.65'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.83%  '
$ws.Range("E49").Value = '  +0.04%  '
$ws.Range("E50").Value = '  -0.34%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '22.89'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -4.49%  '
